$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B to hold "Week_Start_Date".
# This shifts ASIN..is_holiday_week from B..I to C..J.
$ws.Columns.Item(2).Insert()

# New header
$ws.Range("B1").Value = "Week_Start_Date"

# Week start dates for each row (2..17), keep them as plain text (not auto
# converted to Excel date serials) by forcing a text number format first.
$weekStarts = @{
    2  = "2025-01-05"
    3  = "2025-01-12"
    4  = "2025-01-19"
    5  = "2025-01-26"
    6  = "2025-02-02"
    7  = "2025-02-09"
    8  = "2025-02-16"
    9  = "2025-02-23"
    10 = "2025-03-02"
    11 = "2025-03-09"
    12 = "2025-03-16"
    13 = "2025-03-23"
    14 = "2025-03-30"
    15 = "2025-04-06"
    16 = "2025-04-13"
    17 = "2025-04-20"
}

foreach ($r in $weekStarts.Keys) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $weekStarts[$r]
}

# Strip the leading zero from the week labels in column A for weeks 1-9
# (W01 -> W1, ..., W09 -> W9). W10..W16 are already in their final form.
for ($r = 2; $r -le 10; $r++) {
    $week = $r - 1
    $ws.Cells.Item($r, 1).Value = "W$week"
}
